$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper/staging cell used to push numeric-looking strings into cells
# as genuine text (no Excel auto-number coercion) without altering the
# destination cell style: format the helper as Text, paste-special just
# the values into the target, then clear the helper again.
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

$ws.Range("D2").Value = '51.898.45'
$ws.Range("E2").Value = '  -0.35%  '
$ws.Range("D3").Value = '2.784.63'
$ws.Range("E3").Value = '  -2.14%  '
$ws.Range("E4").Value = '  +0.02%  '
$helper.Value = '357.40'
$helper.Copy()
$ws.Range("D5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E5").Value = '  -1.60%  '
$helper.Value = '109.37'
$helper.Copy()
$ws.Range("D6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E6").Value = '  -3.28%  '
$helper.Value = '0.559'
$helper.Copy()
$ws.Range("D7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E7").Value = '  -2.25%  '
$helper.Value = '1.00'
$helper.Copy()
$ws.Range("D8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E8").Value = '  -0.03%  '
$helper.Value = '0.589'
$helper.Copy()
$ws.Range("D9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E9").Value = '  -2.30%  '
$helper.Value = '40.36'
$helper.Copy()
$ws.Range("D10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E10").Value = '  -2.94%  '
$ws.Range("E11").Value = '  -1.96%  '
$ws.Range("E12").Value = '  +1.42%  '
$ws.Range("E13").Value = '  -3.39%  '
$ws.Range("E14").Value = '  -3.13%  '
$ws.Range("D15").Value = '3.225.40'
$ws.Range("E15").Value = '  -1.98%  '
$ws.Range("D16").Value = '2.778.94'
$ws.Range("E16").Value = '  -2.28%  '
$helper.Value = '0.945'
$helper.Copy()
$ws.Range("D17").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E17").Value = '  +3.09%  '
$ws.Range("D18").Value = '51.827.95'
$ws.Range("E18").Value = '  -0.62%  '
$helper.Value = '7.51'
$helper.Copy()
$ws.Range("D19").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E19").Value = '  -0.85%  '
$ws.Range("E20").Value = '  -2.35%  '
$ws.Range("E21").Value = '  -3.75%  '
$ws.Range("D22").Value = '0.0₃0976'
$ws.Range("E22").Value = '  -2.54%  '
$helper.Value = '70.13'
$helper.Copy()
$ws.Range("D23").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E23").Value = '  -0.39%  '
$helper.Value = '269.91'
$helper.Copy()
$ws.Range("D24").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E24").Value = '  +0.18%  '
$ws.Range("E25").Value = '  -4.07%  '
$helper.Value = '26.48'
$helper.Copy()
$ws.Range("D26").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E26").Value = '  -2.47%  '
$ws.Range("E27").Value = '  -0.07%  '
$helper.Value = '0.162'
$helper.Copy()
$ws.Range("D28").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E28").Value = '  +15.44%  '
$ws.Range("E29").Value = '  -0.86%  '
$helper.Value = '2.13'
$helper.Copy()
$ws.Range("D30").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E30").Value = '  -5.27%  '
$ws.Range("E31").Value = '  -4.11%  '
$helper.Value = '52.13'
$helper.Copy()
$ws.Range("D32").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E32").Value = '  -3.21%  '
$helper.Value = '34.56'
$helper.Copy()
$ws.Range("D33").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E33").Value = '  -2.48%  '
$ws.Range("E34").Value = '  -2.69%  '
$helper.Value = '0.0844'
$helper.Copy()
$ws.Range("D35").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E35").Value = '  -0.32%  '
$helper.Value = '5.19'
$helper.Copy()
$ws.Range("D36").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E36").Value = '  -5.61%  '
$ws.Range("E37").Value = '  -0.06%  '
$helper.Value = '18.89'
$helper.Copy()
$ws.Range("D38").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E38").Value = '  +2.02%  '
$ws.Range("E39").Value = '  -2.93%  '
$ws.Range("E40").Value = '  -4.36%  '
$ws.Range("E41").Value = '  +3.13%  '
$ws.Range("E42").Value = '  -2.21%  '
$helper.Value = '2.24'
$helper.Copy()
$ws.Range("D43").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E43").Value = '  -1.43%  '
$helper.Value = '120.05'
$helper.Copy()
$ws.Range("D44").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E44").Value = '  -5.13%  '
$helper.Value = '21.87'
$helper.Copy()
$ws.Range("D45").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E45").Value = '  -8.32%  '
$ws.Range("D46").Value = '2.086.24'
$ws.Range("E46").Value = '  -1.45%  '
$helper.Value = '3.27'
$helper.Copy()
$ws.Range("D47").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E47").Value = '  -4.63%  '
$helper.Value = '0.964'
$helper.Copy()
$ws.Range("D49").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("E49").Value = '  -2.59%  '
$ws.Range("E50").Value = '  -2.75%  '
$ws.Range("E51").Value = '  +31.52%  '

# Clean up the staging cell so it leaves no trace in the sheet.
$helper.Clear()
$excel.CutCopyMode = $false
